$d = $word.ActiveDocument

# --- 1. "within one team" -> "within a team or two" -----------------------
$r = $d.Content
[void]$r.Find.Execute("within one team, it still outperforms", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $r.Start
# layout inside the match: "within "(7) "one"(3) " team"(5) ", it still outperforms"(rest)
# insert first (rightmost edit) so the "one"->"a" shrink (to its left) cannot invalidate its offset
$insPos = $start + 15
$rIns = $d.Range($insPos, $insPos)
[void]$rIns.InsertBefore(" or two")
$rOne = $d.Range($start + 7, $start + 10)
$rOne.Text = "a"

# --- 2. Drop the _GoBack bookmark after "...across the country." ----------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 3. Append the conclusion paragraph text -------------------------------
$n = $d.Paragraphs.Count
$lastReal = $d.Paragraphs.Item($n - 1)
$insertAt = $lastReal.Range.End - 1
$rConclusion = $d.Range($insertAt, $insertAt)
$conclusionText = "Our model, which is a result of a bootstrapped logistical regression using RPI ranking, strength of schedule, and offensive efficiency, has proven to be effective at predicting which teams will enter the NCAA tournament. It is particularly effective at predicting which teams are “on the bubble,” meaning their chances of making the tournament are close to 50/50. While we wish the model was more effective at determining which of those bubble teams would make the tournament, we believe we have made considerable progress in determining the criteria the selection committee uses to make decisions. Some of the variables we initially thought would have strong predictive power actually were not helpful at all in predicting the selection committee’s decisions. Although we may never be able to fully quantify the factors which determine entry into the NCAA tournament, we think our model proves that certain statistics are undeniably more relevant than others and can be used to effectively predict entry into the NCAA tournament. Our final conclusion is that bracketology seems to be a legitimate science and that using statistics to predict the bracket offers a valuable perspective on something countless individuals have and will continue to discuss for years."
[void]$rConclusion.InsertBefore($conclusionText)

# Re-add the _GoBack bookmark at the very end of the document
$docEnd = $d.Content.End
$rEnd = $d.Range($docEnd - 1, $docEnd - 1)
[void]$d.Bookmarks.Add("_GoBack", $rEnd)

# --- 4. Header with right-aligned page number ------------------------------
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$hdr.Range.Paragraphs.Item(1).Style = "Header"
[void]$hdr.PageNumbers.Add(2)
